$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for Table[2] and Table[3], following the same layout as the
# existing Table[1] block (rows 2-4): column B holds the XPath, column C
# holds the validation regex pattern.
$newRows = @(
    @("/NewDataSet/Table[2]/Town", "[A-Z a-z].*"),
    @("/NewDataSet/Table[2]/County", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[2]/PostCode", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[3]/Town", "[A-Z a-z].*"),
    @("/NewDataSet/Table[3]/County", "[A-Z a-z 0-9].*"),
    @("/NewDataSet/Table[3]/PostCode", "[A-Z a-z 0-9].*")
)

$row = 5
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $row = $row + 1
}

# Update the selected range to reflect the newly populated area.
$ws.Range("B5:C10").Select()
